$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the "Conversión del día" report text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 13.04 = 53410.4 pesos`n✅ 53410.4 pesos = 13.07 = 979.3 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update rates in N10, O10, N12, O12 ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 76.68899999999999
$wsTasas.Range("O10").Value = 4095.99
$wsTasas.Range("N12").Value = 4085
$wsTasas.Range("O12").Value = 74.90000000000001
